# Update Name of Algo
# Applies updated RandomForest imputation values to the result_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.10989999999999
$ws.Range("A8").Value = -22.29910000000001
$ws.Range("A10").Value = -21.47909999999997
$ws.Range("A12").Value = -21.56740000000001
$ws.Range("C14").Value = -12.9431
$ws.Range("C15").Value = -14.45429999999998
$ws.Range("A18").Value = -21.72729999999999
$ws.Range("C18").Value = -11.1422
$ws.Range("C20").Value = -11.698
$ws.Range("A25").Value = -21.6027
$ws.Range("C29").Value = -11.6522
$ws.Range("C30").Value = -12.71129999999999
$ws.Range("C31").Value = -12.8541
$ws.Range("C35").Value = -11.66490000000001
$ws.Range("A37").Value = -19.6538
$ws.Range("C40").Value = -13.2993
$ws.Range("C44").Value = -13.15999999999999
$ws.Range("C50").Value = -13.8109
$ws.Range("C54").Value = -13.3315
$ws.Range("A55").Value = -22.2999
$ws.Range("A68").Value = -21.61469999999999
$ws.Range("C68").Value = -11.2929
$ws.Range("C76").Value = -12.5018
$ws.Range("A77").Value = -20.47379999999999
$ws.Range("A78").Value = -20.63299999999998
$ws.Range("A79").Value = -20.59489999999999
$ws.Range("A80").Value = -19.9285
$ws.Range("A81").Value = -21.7421
$ws.Range("A82").Value = -22.20620000000001
$ws.Range("A84").Value = -21.82869999999999
$ws.Range("C87").Value = -13.59229999999999
$ws.Range("C88").Value = -13.38449999999999
$ws.Range("C92").Value = -11.063
$ws.Range("C96").Value = -12.81920000000001
$ws.Range("C98").Value = -12.1404
$ws.Range("A101").Value = -20.64069999999997
$ws.Range("C101").Value = -12.86750000000001
$ws.Range("A102").Value = -19.50199999999999
$ws.Range("C102").Value = -13.2855
